$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 2777.7778
$ws.Range("I64").Value = 2733.3333
$ws.Range("J64").Value = 2866.6667
$ws.Range("K64").Value = 2733.3333
$ws.Range("L64").Value = 2866.6667
$ws.Range("M64").Value = -2485.3333
$ws.Range("N64").Value = -3362.6667
$ws.Range("H67").Value = 2777.7778
$ws.Range("I67").Value = 2733.3333
$ws.Range("J67").Value = 2866.6667
$ws.Range("K67").Value = 2733.3333
$ws.Range("L67").Value = 2866.6667
$ws.Range("M67").Value = -1875.3333
$ws.Range("N67").Value = -4582.6667
$ws.Range("H112").Value = 1632.4579
$ws.Range("I112").Value = 700
$ws.Range("J112").Value = 1643.8292
$ws.Range("K112").Value = 2100
$ws.Range("L112").Value = 4931.487599999999
$ws.Range("M112").Value = -992
$ws.Range("N112").Value = -7147.487599999999
$ws.Range("H121").Value = 3700.6155
$ws.Range("I121").Value = 0
$ws.Range("J121").Value = 3700.6155
$ws.Range("K121").Value = 0
$ws.Range("L121").Value = 11101.8465
$ws.Range("N121").Value = -14595.8465

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 698.5789
$ws.Range("I2").Value = 640
$ws.Range("J2").Value = 763.6667
$ws.Range("K2").Value = 640
$ws.Range("L2").Value = 763.6667
$ws.Range("M2").Value = -527
$ws.Range("N2").Value = -989.6667
$ws.Range("H116").Value = 698.5789
$ws.Range("I116").Value = 640
$ws.Range("J116").Value = 763.6667
$ws.Range("K116").Value = 640
$ws.Range("L116").Value = 763.6667
$ws.Range("M116").Value = 1654
$ws.Range("N116").Value = -5351.6667
$ws.Range("H137").Value = 41045.6
$ws.Range("I137").Value = 0
$ws.Range("J137").Value = 41045.6
$ws.Range("K137").Value = 0
$ws.Range("L137").Value = 41045.6
$ws.Range("N137").Value = -51245.6

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 698.5789
$ws.Range("I3").Value = 640
$ws.Range("J3").Value = 763.6667
$ws.Range("K3").Value = 640
$ws.Range("L3").Value = 763.6667
$ws.Range("M3").Value = -526
$ws.Range("N3").Value = -991.6667
$ws.Range("H22").Value = 1501.5
$ws.Range("I22").Value = 1001
$ws.Range("J22").Value = 2002
$ws.Range("K22").Value = 1001
$ws.Range("L22").Value = 2002
$ws.Range("M22").Value = -828
$ws.Range("N22").Value = -2348
$ws.Range("H134").Value = 2693.8708
$ws.Range("I134").Value = 1592.8235
$ws.Range("J134").Value = 4030.8572
$ws.Range("K134").Value = 4778.470499999999
$ws.Range("L134").Value = 12092.5716
$ws.Range("M134").Value = -2243.470499999999
$ws.Range("N134").Value = -17162.5716

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 4039.125
$ws.Range("I132").Value = 2267
$ws.Range("J132").Value = 4925.1875
$ws.Range("K132").Value = 6801
$ws.Range("L132").Value = 14775.5625
$ws.Range("M132").Value = -4271
$ws.Range("N132").Value = -19835.5625
$ws.Range("H134").Value = 7095.8184
$ws.Range("I134").Value = 9376.5
$ws.Range("J134").Value = 4359
$ws.Range("K134").Value = 28129.5
$ws.Range("L134").Value = 13077
$ws.Range("M134").Value = -25594.5
$ws.Range("N134").Value = -18147

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1794.8966
$ws.Range("I5").Value = 335.4375
$ws.Range("J5").Value = 3591.1538
$ws.Range("K5").Value = 1006.3125
$ws.Range("L5").Value = 10773.4614
$ws.Range("M5").Value = -894.3125
$ws.Range("N5").Value = -10997.4614
$ws.Range("H107").Value = 46542.59
$ws.Range("I107").Value = 440.25
$ws.Range("J107").Value = 169482.17
$ws.Range("K107").Value = 1320.75
$ws.Range("L107").Value = 508446.51
$ws.Range("M107").Value = 599.25
$ws.Range("N107").Value = -512286.51
$ws.Range("H113").Value = 567.2
$ws.Range("I113").Value = 547.9259
$ws.Range("J113").Value = 632.25
$ws.Range("K113").Value = 1643.7777
$ws.Range("L113").Value = 1896.75
$ws.Range("M113").Value = 526.2223000000001
$ws.Range("N113").Value = -6236.75
$ws.Range("H131").Value = 1035.3889
$ws.Range("I131").Value = 2787.5
$ws.Range("J131").Value = 895.22
$ws.Range("K131").Value = 8362.5
$ws.Range("L131").Value = 2685.66
$ws.Range("M131").Value = -3322.5
$ws.Range("N131").Value = -12765.66
$ws.Range("H135").Value = 1794.8966
$ws.Range("I135").Value = 335.4375
$ws.Range("J135").Value = 3591.1538
$ws.Range("K135").Value = 3018.9375
$ws.Range("L135").Value = 32320.3842
$ws.Range("M135").Value = -483.9375
$ws.Range("N135").Value = -37390.3842

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H134").Value = 44888.09
$ws.Range("I134").Value = 19296
$ws.Range("J134").Value = 46106.76
$ws.Range("K134").Value = 57888
$ws.Range("L134").Value = 138320.28
$ws.Range("M134").Value = -55353
$ws.Range("N134").Value = -143390.28
$ws.Range("H137").Value = 59433.11
$ws.Range("I137").Value = 48000
$ws.Range("J137").Value = 60105.65
$ws.Range("K137").Value = 48000
$ws.Range("L137").Value = 60105.65
$ws.Range("M137").Value = -42900
$ws.Range("N137").Value = -70305.64999999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H23").Value = 12000
$ws.Range("I23").Value = 0
$ws.Range("J23").Value = 12000
$ws.Range("K23").Value = 0
$ws.Range("L23").Value = 12000
$ws.Range("N23").Value = -12460
$ws.Range("H40").Value = 8178.05
$ws.Range("I40").Value = 6958.7144
$ws.Range("J40").Value = 8834.615
$ws.Range("K40").Value = 6958.7144
$ws.Range("L40").Value = 8834.615
$ws.Range("M40").Value = -6822.7144
$ws.Range("N40").Value = -9106.615
$ws.Range("H122").Value = 5835.357
$ws.Range("I122").Value = 3385.7144
$ws.Range("J122").Value = 8285
$ws.Range("K122").Value = 10157.1432
$ws.Range("L122").Value = 24855
$ws.Range("M122").Value = -7707.143199999999
$ws.Range("N122").Value = -29755
$ws.Range("H132").Value = 2994.5615
$ws.Range("I132").Value = 1577.5526
$ws.Range("J132").Value = 4533.029
$ws.Range("K132").Value = 4732.6578
$ws.Range("L132").Value = 13599.087
$ws.Range("M132").Value = -2202.6578
$ws.Range("N132").Value = -18659.087

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H39").Value = 24999.5
$ws.Range("I39").Value = 0
$ws.Range("J39").Value = 24999.5
$ws.Range("K39").Value = 0
$ws.Range("L39").Value = 24999.5
$ws.Range("N39").Value = -25825.5
$ws.Range("H80").Value = 44559.8
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 44559.8
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 44559.8
$ws.Range("N80").Value = -46555.8
$ws.Range("H83").Value = 44559.8
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 44559.8
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 133679.4
$ws.Range("N83").Value = -143663.4
$ws.Range("H108").Value = 39500
$ws.Range("I108").Value = 0
$ws.Range("J108").Value = 39500
$ws.Range("K108").Value = 0
$ws.Range("L108").Value = 39500
$ws.Range("N108").Value = -47180
$ws.Range("H122").Value = 3042.7568
$ws.Range("I122").Value = 1953.1538
$ws.Range("J122").Value = 5618.1816
$ws.Range("K122").Value = 5859.4614
$ws.Range("L122").Value = 16854.5448
$ws.Range("M122").Value = -3409.4614
$ws.Range("N122").Value = -21754.5448
$ws.Range("H132").Value = 16669622
$ws.Range("I132").Value = 1633.6666
$ws.Range("J132").Value = 23813046
$ws.Range("K132").Value = 4900.9998
$ws.Range("L132").Value = 71439138
$ws.Range("M132").Value = -2370.9998
$ws.Range("N132").Value = -71444198
$ws.Range("H136").Value = 1466.2222
$ws.Range("I136").Value = 652.2222
$ws.Range("J136").Value = 3094.2222
$ws.Range("K136").Value = 1956.6666
$ws.Range("L136").Value = 9282.6666
$ws.Range("M136").Value = 593.3334
$ws.Range("N136").Value = -14382.6666
$ws.Range("M39").ClearContents()

Write-Host "Applied all cell updates"